# "contingencies with rene fine"
# Lay out a tiny 2x2 block on Sheet1:
#   B1 = 0                         (bold, thin box border, centered/top)
#   A2 = 0                         (same style as B1)
#   B2 = "disconnected_elements"   (plain text, default style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / thin-bordered / centered+top style on B1 first.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.Borders.LineStyle = 1        # xlContinuous
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop

# Re-use the exact same style on A2 (copy/paste-format keeps both cells
# pointing at the same cellXfs entry instead of minting a duplicate one).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
